$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Seq1 (B16) data
$ws.Range("C16").Formula = "=COUNT(D16:CY16)"
$rng = $ws.Range("D16:AO16")
$rng.NumberFormat = "h:mm"
$arr = New-Object 'object[,]' 1,38
$arr[0,0] = [double]"1.1111111111111112E-2"
$arr[0,1] = [double]"4.7222222222222221E-2"
$arr[0,2] = [double]"5.486111111111111E-2"
$arr[0,3] = [double]"6.25E-2"
$arr[0,4] = [double]"8.3333333333333329E-2"
$arr[0,5] = [double]"9.0277777777777776E-2"
$arr[0,6] = [double]"0.10138888888888889"
$arr[0,7] = [double]"0.11944444444444445"
$arr[0,8] = [double]"0.15138888888888888"
$arr[0,9] = [double]"0.17222222222222222"
$arr[0,10] = [double]"0.20833333333333334"
$arr[0,11] = [double]"0.23680555555555555"
$arr[0,12] = [double]"0.24444444444444444"
$arr[0,13] = [double]"0.26041666666666669"
$arr[0,14] = [double]"0.27708333333333335"
$arr[0,15] = [double]"0.28680555555555554"
$arr[0,16] = [double]"0.29444444444444445"
$arr[0,17] = [double]"0.30555555555555558"
$arr[0,18] = [double]"0.31319444444444444"
$arr[0,19] = [double]"0.32013888888888886"
$arr[0,20] = [double]"0.33263888888888887"
$arr[0,21] = [double]"0.35208333333333336"
$arr[0,22] = [double]"0.36527777777777776"
$arr[0,23] = [double]"0.37291666666666667"
$arr[0,24] = [double]"0.40208333333333335"
$arr[0,25] = [double]"0.40972222222222221"
$arr[0,26] = [double]"0.42638888888888887"
$arr[0,27] = [double]"0.43402777777777779"
$arr[0,28] = [double]"0.46041666666666664"
$arr[0,29] = [double]"0.46805555555555556"
$arr[0,30] = [double]"0.47638888888888886"
$arr[0,31] = [double]"0.50277777777777777"
$arr[0,32] = [double]"0.51041666666666663"
$arr[0,33] = [double]"0.52638888888888891"
$arr[0,34] = [double]"0.53680555555555554"
$arr[0,35] = [double]"0.5444444444444444"
$arr[0,36] = [double]"0.57291666666666663"
$arr[0,37] = [double]"0.58611111111111114"
$rng.Value = $arr

# Rows 17-22: shared COUNT formula, then data for 17-19
$ws.Range("C17:C22").Formula = "=COUNT(D17:CY17)"

# Row 17 data
$rng = $ws.Range("D17:R17")
$rng.NumberFormat = "h:mm"
$arr = New-Object 'object[,]' 1,15
$arr[0,0] = [double]"5.8333333333333334E-2"
$arr[0,1] = [double]"6.805555555555555E-2"
$arr[0,2] = [double]"0.11458333333333333"
$arr[0,3] = [double]"0.17291666666666666"
$arr[0,4] = [double]"0.2048611111111111"
$arr[0,5] = [double]"0.21180555555555555"
$arr[0,6] = [double]"0.22500000000000001"
$arr[0,7] = [double]"0.23958333333333334"
$arr[0,8] = [double]"0.25972222222222224"
$arr[0,9] = [double]"0.27430555555555558"
$arr[0,10] = [double]"0.33055555555555555"
$arr[0,11] = [double]"0.36944444444444446"
$arr[0,12] = [double]"0.42222222222222222"
$arr[0,13] = [double]"0.43055555555555558"
$arr[0,14] = [double]"0.44444444444444442"
$rng.Value = $arr

# Row 18 data
$rng = $ws.Range("D18:BJ18")
$rng.NumberFormat = "h:mm"
$arr = New-Object 'object[,]' 1,59
$arr[0,0] = [double]"1.0416666666666666E-2"
$arr[0,1] = [double]"1.7361111111111112E-2"
$arr[0,2] = [double]"2.5694444444444443E-2"
$arr[0,3] = [double]"3.3333333333333333E-2"
$arr[0,4] = [double]"4.5138888888888888E-2"
$arr[0,5] = [double]"5.2777777777777778E-2"
$arr[0,6] = [double]"6.0416666666666667E-2"
$arr[0,7] = [double]"7.0833333333333331E-2"
$arr[0,8] = [double]"8.5416666666666669E-2"
$arr[0,9] = [double]"9.4444444444444442E-2"
$arr[0,10] = [double]"0.10486111111111111"
$arr[0,11] = [double]"0.12916666666666668"
$arr[0,12] = [double]"0.14097222222222222"
$arr[0,13] = [double]"0.14861111111111111"
$arr[0,14] = [double]"0.17291666666666666"
$arr[0,15] = [double]"0.18055555555555555"
$arr[0,16] = [double]"0.18888888888888888"
$arr[0,17] = [double]"0.20902777777777778"
$arr[0,18] = [double]"0.21666666666666667"
$arr[0,19] = [double]"0.22638888888888889"
$arr[0,20] = [double]"0.23402777777777778"
$arr[0,21] = [double]"0.25277777777777777"
$arr[0,22] = [double]"0.26250000000000001"
$arr[0,23] = [double]"0.27083333333333331"
$arr[0,24] = [double]"0.27847222222222223"
$arr[0,25] = [double]"0.30208333333333331"
$arr[0,26] = [double]"0.30972222222222223"
$arr[0,27] = [double]"0.31736111111111109"
$arr[0,28] = [double]"0.32847222222222222"
$arr[0,29] = [double]"0.33611111111111114"
$arr[0,30] = [double]"0.35208333333333336"
$arr[0,31] = [double]"0.38333333333333336"
$arr[0,32] = [double]"0.41875000000000001"
$arr[0,33] = [double]"0.42638888888888887"
$arr[0,34] = [double]"0.44027777777777777"
$arr[0,35] = [double]"0.44722222222222224"
$arr[0,36] = [double]"0.47013888888888888"
$arr[0,37] = [double]"0.4826388888888889"
$arr[0,38] = [double]"0.49027777777777776"
$arr[0,39] = [double]"0.49791666666666667"
$arr[0,40] = [double]"0.50694444444444442"
$arr[0,41] = [double]"0.51458333333333328"
$arr[0,42] = [double]"0.52222222222222225"
$arr[0,43] = [double]"0.54722222222222228"
$arr[0,44] = [double]"0.55486111111111114"
$arr[0,45] = [double]"0.56874999999999998"
$arr[0,46] = [double]"0.57638888888888884"
$arr[0,47] = [double]"0.58402777777777781"
$arr[0,48] = [double]"0.59305555555555556"
$arr[0,49] = [double]"0.61319444444444449"
$arr[0,50] = [double]"0.62083333333333335"
$arr[0,51] = [double]"0.62847222222222221"
$arr[0,52] = [double]"0.64236111111111116"
$arr[0,53] = [double]"0.64930555555555558"
$arr[0,54] = [double]"0.65694444444444444"
$arr[0,55] = [double]"0.6645833333333333"
$arr[0,56] = [double]"0.67222222222222228"
$arr[0,57] = [double]"0.68125000000000002"
$arr[0,58] = [double]"0.70416666666666672"
$rng.Value = $arr

# Row 19 data
$rng = $ws.Range("D19:BJ19")
$rng.NumberFormat = "h:mm"
$arr = New-Object 'object[,]' 1,59
$arr[0,0] = [double]"5.5555555555555558E-3"
$arr[0,1] = [double]"1.3194444444444444E-2"
$arr[0,2] = [double]"2.7083333333333334E-2"
$arr[0,3] = [double]"3.4722222222222224E-2"
$arr[0,4] = [double]"5.1388888888888887E-2"
$arr[0,5] = [double]"5.9027777777777776E-2"
$arr[0,6] = [double]"7.1527777777777773E-2"
$arr[0,7] = [double]"8.4027777777777785E-2"
$arr[0,8] = [double]"0.10277777777777777"
$arr[0,9] = [double]"0.11041666666666666"
$arr[0,10] = [double]"0.12291666666666666"
$arr[0,11] = [double]"0.13263888888888889"
$arr[0,12] = [double]"0.14027777777777778"
$arr[0,13] = [double]"0.14791666666666667"
$arr[0,14] = [double]"0.15555555555555556"
$arr[0,15] = [double]"0.17291666666666666"
$arr[0,16] = [double]"0.18055555555555555"
$arr[0,17] = [double]"0.18819444444444444"
$arr[0,18] = [double]"0.19583333333333333"
$arr[0,19] = [double]"0.21527777777777779"
$arr[0,20] = [double]"0.23194444444444445"
$arr[0,21] = [double]"0.23958333333333334"
$arr[0,22] = [double]"0.24722222222222223"
$arr[0,23] = [double]"0.25486111111111109"
$arr[0,24] = [double]"0.26250000000000001"
$arr[0,25] = [double]"0.27013888888888887"
$arr[0,26] = [double]"0.28194444444444444"
$arr[0,27] = [double]"0.2902777777777778"
$arr[0,28] = [double]"0.29791666666666666"
$arr[0,29] = [double]"0.30555555555555558"
$arr[0,30] = [double]"0.33263888888888887"
$arr[0,31] = [double]"0.34027777777777779"
$arr[0,32] = [double]"0.34791666666666665"
$arr[0,33] = [double]"0.35555555555555557"
$arr[0,34] = [double]"0.3659722222222222"
$arr[0,35] = [double]"0.37777777777777777"
$arr[0,36] = [double]"0.38541666666666669"
$arr[0,37] = [double]"0.39305555555555555"
$arr[0,38] = [double]"0.40069444444444446"
$arr[0,39] = [double]"0.41805555555555557"
$arr[0,40] = [double]"0.42569444444444443"
$arr[0,41] = [double]"0.43680555555555556"
$arr[0,42] = [double]"0.44722222222222224"
$arr[0,43] = [double]"0.4548611111111111"
$arr[0,44] = [double]"0.46250000000000002"
$arr[0,45] = [double]"0.47013888888888888"
$arr[0,46] = [double]"0.4777777777777778"
$arr[0,47] = [double]"0.49444444444444446"
$arr[0,48] = [double]"0.50208333333333333"
$arr[0,49] = [double]"0.50972222222222219"
$arr[0,50] = [double]"0.51736111111111116"
$arr[0,51] = [double]"0.53472222222222221"
$arr[0,52] = [double]"0.54513888888888884"
$arr[0,53] = [double]"0.55277777777777781"
$arr[0,54] = [double]"0.56041666666666667"
$arr[0,55] = [double]"0.58125000000000004"
$arr[0,56] = [double]"0.59375"
$arr[0,57] = [double]"0.6069444444444444"
$arr[0,58] = [double]"0.61597222222222225"
$rng.Value = $arr

# Row 23: total
$ws.Range("C23").Formula = "=SUM(C16:C22)"

# Hyperlink for B16
$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.youtube.com/watch?v=awtmTJW9ic8")
$ws.Range("B16").Style = "Hyperlink"

# Update selection to C13 (also updates sheet view)
$ws.Range("C13").Select()

